$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new "class_numeric" column (header) after the existing "roll" column
$ws.Range("H1").Value = "class_numeric"

# Match the new best-fit width Excel computed for the freshly added column
$ws.Columns.Item(8).AutoFit()

# Reflect the new active cell/selection left behind after the edit
$ws.Range("H3").Select()
